$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 27064
$ws.Range("B2").Value = "Luana Cavalcanti"
$ws.Range("C2").Value = "Vendas"
$ws.Range("D2").Value = "Doenca"
$ws.Range("E2").Value = 6
$ws.Range("F2").Value = 45106
$ws.Range("G2").Value = 8368.51

# Row 3
$ws.Range("A3").Value = 2754
$ws.Range("B3").Value = "Srta. Ana Cecília Marques"
$ws.Range("C3").Value = "Operacoes"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 45090
$ws.Range("G3").Value = 4309.15

# Row 4
$ws.Range("A4").Value = 41226
$ws.Range("B4").Value = "Luan Alves"
$ws.Range("C4").Value = "Engenharia"
$ws.Range("E4").Value = 4
$ws.Range("F4").Value = 45085
$ws.Range("G4").Value = 4319.28

# Row 5
$ws.Range("A5").Value = 28251
$ws.Range("B5").Value = "Juliana Marques"
$ws.Range("C5").Value = "Operacoes"
$ws.Range("D5").Value = "Consulta medica"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 45086
$ws.Range("G5").Value = 6473.13

# Row 6
$ws.Range("A6").Value = 41390
$ws.Range("B6").Value = "Enzo Gabriel Garcia"
$ws.Range("C6").Value = "TI"
$ws.Range("D6").Value = "Doenca"
$ws.Range("E6").Value = 7
$ws.Range("F6").Value = 45106
$ws.Range("G6").Value = 7022.3

# Row 7
$ws.Range("A7").Value = 57562
$ws.Range("B7").Value = "Dra. Julia Rezende"
$ws.Range("C7").Value = "Engenharia"
$ws.Range("D7").Value = "Viagem de negocios"
$ws.Range("E7").Value = 5
$ws.Range("F7").Value = 45084
$ws.Range("G7").Value = 8096.78

# Row 8
$ws.Range("A8").Value = 80684
$ws.Range("B8").Value = "Srta. Stephany Araújo"
$ws.Range("E8").Value = 5
$ws.Range("F8").Value = 45084
$ws.Range("G8").Value = 8010.75

# Row 9
$ws.Range("A9").Value = 15129
$ws.Range("B9").Value = "Bella Costela"
$ws.Range("C9").Value = "Atendimento ao Cliente"
$ws.Range("D9").Value = "Problemas pessoais"
$ws.Range("E9").Value = 7
$ws.Range("F9").Value = 45088
$ws.Range("G9").Value = 2685.75

# Row 10
$ws.Range("A10").Value = 5861
$ws.Range("B10").Value = "Gabriel Leão"
$ws.Range("C10").Value = "TI"
$ws.Range("D10").Value = "Consulta medica"
$ws.Range("E10").Value = 4
$ws.Range("F10").Value = 45094
$ws.Range("G10").Value = 6576.09

# Row 11
$ws.Range("A11").Value = 82266
$ws.Range("B11").Value = "Ana Luiza da Conceição"
$ws.Range("C11").Value = "Recursos Humanos"
$ws.Range("D11").Value = "Consulta medica"
$ws.Range("F11").Value = 45105
$ws.Range("G11").Value = 6344.75
